$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40-46 down to 41-47
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new weekly price record
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 44551
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100103
$ws.Cells.Item(40, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value = 100103003
$ws.Cells.Item(40, 10).Value = "Damasco"
$ws.Cells.Item(40, 11).Value = "Castle Brite"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 45
$ws.Cells.Item(40, 14).Value = 18000
$ws.Cells.Item(40, 15).Value = 18000
$ws.Cells.Item(40, 16).Value = 18000
$ws.Cells.Item(40, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(40, 19).Value = 1000
$ws.Cells.Item(40, 20).Value = 18
